$wb = $excel.ActiveWorkbook

$ws = $wb.Sheets.Item("ALC")
# Row 18
$ws.Range("H18").Value = 831.0526
$ws.Range("I18").Value = 954
$ws.Range("J18").Value = 370
$ws.Range("K18").Value = 954
$ws.Range("L18").Value = 370
$ws.Range("M18").Value = -670
$ws.Range("N18").Value = -938

# Row 33
$ws.Range("H33").Value = 1541.9524
$ws.Range("I33").Value = 1010.4706
$ws.Range("K33").Value = 1010.4706
$ws.Range("M33").Value = -781.4706

# Row 86
$ws.Range("H86").Value = 6675187
$ws.Range("I86").Value = 11906480
$ws.Range("J86").Value = 1443894.4
$ws.Range("K86").Value = 11906480
$ws.Range("L86").Value = 1443894.4
$ws.Range("M86").Value = -11905357
$ws.Range("N86").Value = -1446140.4

# Row 89
$ws.Range("H89").Value = 6675187
$ws.Range("I89").Value = 11906480
$ws.Range("J89").Value = 1443894.4
$ws.Range("K89").Value = 59532400
$ws.Range("L89").Value = 7219472
$ws.Range("M89").Value = -59526784
$ws.Range("N89").Value = -7230704

# Row 94
$ws.Range("H94").Value = 1022
$ws.Range("I94").Value = 1022
$ws.Range("K94").Value = 1022
$ws.Range("M94").Value = -571

# Row 98
$ws.Range("H98").Value = 25411528
$ws.Range("I98").Value = 9091964
$ws.Range("J98").Value = 85249930
$ws.Range("K98").Value = 9091964
$ws.Range("L98").Value = 85249930
$ws.Range("M98").Value = -9090466
$ws.Range("N98").Value = -85252926

# Row 106
$ws.Range("H106").Value = 43480584
$ws.Range("I106").Value = 62501780
$ws.Range("J106").Value = 3556.5715
$ws.Range("K106").Value = 62501780
$ws.Range("L106").Value = 3556.5715
$ws.Range("M106").Value = -62501149
$ws.Range("N106").Value = -4818.5715

# Row 122
$ws.Range("H122").Value = 25411528
$ws.Range("I122").Value = 9091964
$ws.Range("J122").Value = 85249930
$ws.Range("K122").Value = 27275892
$ws.Range("L122").Value = 255749790
$ws.Range("M122").Value = -27273442
$ws.Range("N122").Value = -255754690

# Row 138
$ws.Range("H138").Value = 2127.3552
$ws.Range("I138").Value = 1453.64
$ws.Range("J138").Value = 3422.9614
$ws.Range("K138").Value = 4360.92
$ws.Range("L138").Value = 10268.8842
$ws.Range("M138").Value = 779.0799999999999
$ws.Range("N138").Value = -20548.8842

$ws = $wb.Sheets.Item("ARM")
# Row 4
$ws.Range("H4").Value = 900
$ws.Range("I4").Value = 900
$ws.Range("J4").Value = 900
$ws.Range("K4").Value = 900
$ws.Range("L4").Value = 900
$ws.Range("M4").Value = -784
$ws.Range("N4").Value = -1132

# Row 5
$ws.Range("H5").Value = 389.9091
$ws.Range("I5").Value = 197.8
$ws.Range("J5").Value = 550
$ws.Range("K5").Value = 197.8
$ws.Range("L5").Value = 550
$ws.Range("M5").Value = -85.80000000000001
$ws.Range("N5").Value = -774

# Row 32
$ws.Range("H32").Value = 723.6799999999999
$ws.Range("I32").Value = 714.9796
$ws.Range("J32").Value = 1150
$ws.Range("K32").Value = 714.9796
$ws.Range("L32").Value = 1150
$ws.Range("M32").Value = -427.9796
$ws.Range("N32").Value = -1724

# Row 34
$ws.Range("H34").Value = 19800
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 19800
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 19800
$ws.Range("M34").ClearContents()
$ws.Range("N34").Value = -20342

# Row 45
$ws.Range("H45").Value = 1001234.4
$ws.Range("I45").Value = 1251141.1
$ws.Range("J45").Value = 1607.5
$ws.Range("K45").Value = 1251141.1
$ws.Range("L45").Value = 1607.5
$ws.Range("M45").Value = -1250764.1
$ws.Range("N45").Value = -2361.5

# Row 61
$ws.Range("H61").Value = 2393504
$ws.Range("I61").Value = 1126879
$ws.Range("J61").Value = 11766528
$ws.Range("K61").Value = 1126879
$ws.Range("L61").Value = 11766528
$ws.Range("M61").Value = -1126667
$ws.Range("N61").Value = -11766952

# Row 110
$ws.Range("H110").Value = 1882.3572
$ws.Range("I110").Value = 974
$ws.Range("J110").Value = 4153.25
$ws.Range("K110").Value = 974
$ws.Range("L110").Value = 4153.25
$ws.Range("M110").Value = 1071
$ws.Range("N110").Value = -8243.25

# Row 136
$ws.Range("H136").Value = 2393504
$ws.Range("I136").Value = 1126879
$ws.Range("J136").Value = 11766528
$ws.Range("K136").Value = 3380637
$ws.Range("L136").Value = 35299584
$ws.Range("M136").Value = -3378087
$ws.Range("N136").Value = -35304684

$ws = $wb.Sheets.Item("BSM")
# Row 4
$ws.Range("H4").Value = 389.9091
$ws.Range("I4").Value = 197.8
$ws.Range("J4").Value = 550
$ws.Range("K4").Value = 197.8
$ws.Range("L4").Value = 550
$ws.Range("M4").Value = -82.80000000000001
$ws.Range("N4").Value = -780

# Row 22
$ws.Range("H22").Value = 509.13043
$ws.Range("I22").Value = 508.36365
$ws.Range("J22").Value = 526
$ws.Range("K22").Value = 508.36365
$ws.Range("L22").Value = 526
$ws.Range("M22").Value = -335.36365
$ws.Range("N22").Value = -872

$ws = $wb.Sheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 669.6799999999999
$ws.Range("I16").Value = 652
$ws.Range("J16").Value = 762.5
$ws.Range("K16").Value = 652
$ws.Range("L16").Value = 762.5
$ws.Range("M16").Value = -365
$ws.Range("N16").Value = -1336.5

# Row 22
$ws.Range("H22").Value = 163.72728
$ws.Range("I22").Value = 166.66667
$ws.Range("J22").Value = 150.5
$ws.Range("K22").Value = 166.66667
$ws.Range("L22").Value = 150.5
$ws.Range("M22").Value = 183.33333
$ws.Range("N22").Value = -850.5

# Row 113
$ws.Range("H113").Value = 669.6799999999999
$ws.Range("I113").Value = 652
$ws.Range("J113").Value = 762.5
$ws.Range("K113").Value = 652
$ws.Range("L113").Value = 762.5
$ws.Range("M113").Value = 1518
$ws.Range("N113").Value = -5102.5

$ws = $wb.Sheets.Item("GSM")
# Row 2
$ws.Range("H2").Value = 103.666664
$ws.Range("I2").Value = 55.5
$ws.Range("J2").Value = 200
$ws.Range("K2").Value = 55.5
$ws.Range("L2").Value = 200
$ws.Range("M2").Value = 57.5
$ws.Range("N2").Value = -426

$ws = $wb.Sheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 4695.5557
$ws.Range("I22").Value = 5875
$ws.Range("J22").Value = 4358.5713
$ws.Range("K22").Value = 5875
$ws.Range("L22").Value = 4358.5713
$ws.Range("M22").Value = -5580
$ws.Range("N22").Value = -4948.5713

# Row 27
$ws.Range("H27").Value = 4695.5557
$ws.Range("I27").Value = 5875
$ws.Range("J27").Value = 4358.5713
$ws.Range("K27").Value = 5875
$ws.Range("L27").Value = 4358.5713
$ws.Range("M27").Value = -5768
$ws.Range("N27").Value = -4572.5713

# Row 61
$ws.Range("H61").Value = 1449.2858
$ws.Range("I61").Value = 1170
$ws.Range("J61").Value = 1728.5714
$ws.Range("K61").Value = 1170
$ws.Range("L61").Value = 1728.5714
$ws.Range("M61").Value = -968
$ws.Range("N61").Value = -2132.5714

# Row 113
$ws.Range("H113").Value = 1449.2858
$ws.Range("I113").Value = 1170
$ws.Range("J113").Value = 1728.5714
$ws.Range("K113").Value = 1170
$ws.Range("L113").Value = 1728.5714
$ws.Range("M113").Value = 1000
$ws.Range("N113").Value = -6068.5714

$ws = $wb.Sheets.Item("WVR")
# Row 107
$ws.Range("H107").Value = 16213.685
$ws.Range("I107").Value = 25045.25
$ws.Range("J107").Value = 9790.727999999999
$ws.Range("K107").Value = 75135.75
$ws.Range("L107").Value = 29372.184
$ws.Range("M107").Value = -73215.75
$ws.Range("N107").Value = -33212.18399999999

Write-Host "All edits applied successfully"
